$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Bandit"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 8

$ws.Range("F17").Select()
